$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = "TAB-1"
$ws.Range("B2").Value = "https://palegoldenrod-ant-677872.hostingersite.com/basic_page/applied-biosystems-model-7000-real-time-thermalcycler/"
$ws.Range("C2").Value = "https://palegoldenrod-ant-677872.hostingersite.com/chemistry/applied-biosystems-model-7000-real-time-thermalcycler/"
$ws.Range("D2").Value = $true
$ws.Range("E2").Value = "https://palegoldenrod-ant-677872.hostingersite.com/chemistry/applied-biosystems-model-7000-real-time-thermalcycler/"
$ws.Range("H2").Value = 403

# --- Row 3 ---
$ws.Range("A3").Value = "TAB-2"
$ws.Range("B3").Value = "https://palegoldenrod-ant-677872.hostingersite.com/basic_page/105-description/"
$ws.Range("C3").Value = "https://palegoldenrod-ant-677872.hostingersite.com/chemstockroom/105-description/"
$ws.Range("D3").Value = $true
$ws.Range("E3").Value = "https://palegoldenrod-ant-677872.hostingersite.com/chemstockroom/105-description/"
$ws.Range("H3").Value = 403

# --- Row 4 ---
$ws.Range("A4").Value = "TAB-3"
$ws.Range("B4").Value = "https://palegoldenrod-ant-677872.hostingersite.com/basic_page/119-description/"
$ws.Range("C4").Value = "https://palegoldenrod-ant-677872.hostingersite.com/chemstockroom/119-description/"
$ws.Range("D4").Value = $true
$ws.Range("E4").Value = "https://palegoldenrod-ant-677872.hostingersite.com/chemstockroom/119-description/"
$ws.Range("H4").Value = 403

# --- Row 5: rebuilt (old academics/facebook row -> new TAB-5 error row) ---
$ws.Range("A5").Value = "TAB-5"
$ws.Range("B5").Value = "https://palegoldenrod-ant-677872.hostingersite.com/basic_page/216-description/"
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("E5").Value = "https://palegoldenrod-ant-677872.hostingersite.com/basic_page/216-description/"
$ws.Range("F5").Value = $true
$ws.Range("H5").ClearContents()
$esc = [char]27
$errText5 = "page.goto: net::ERR_ABORTED; maybe frame was detached?`nCall log:`n$esc[2m  - navigating to `"https://palegoldenrod-ant-677872.hostingersite.com/basic_page/216-description/`", waiting until `"domcontentloaded`"$esc[22m`n"
$ws.Range("I5").Value = $errText5

# --- Row 6: rebuilt (old academics/linkedin row -> new TAB-4 error row) ---
$ws.Range("A6").Value = "TAB-4"
$ws.Range("B6").Value = "https://palegoldenrod-ant-677872.hostingersite.com/basic_page/209-description/"
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("E6").Value = "https://palegoldenrod-ant-677872.hostingersite.com/basic_page/209-description/"
$ws.Range("F6").Value = $true
$ws.Range("H6").ClearContents()
$errText6 = "page.goto: net::ERR_ABORTED; maybe frame was detached?`nCall log:`n$esc[2m  - navigating to `"https://palegoldenrod-ant-677872.hostingersite.com/basic_page/209-description/`", waiting until `"domcontentloaded`"$esc[22m`n"
$ws.Range("I6").Value = $errText6

# Reset the auto row-height bump triggered by the multi-line error text so the
# rows keep a plain <row> element (no ht/customHeight attributes), matching
# the original sheet's styling.
$ws.Rows(5).EntireRow.AutoFit()
$ws.Rows(6).EntireRow.AutoFit()

# --- Remove rows 7-9 entirely ---
$ws.Rows("7:9").Delete()

Write-Host "done"
